$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- D5: update numeric value ---
$ws.Range("D5").Value = 0.4609

# --- G15: fix double-dot -> triple-dot in voltage range text ---
$g15 = "{'Тип диода': 'LED', 'Диаметр диода LED': '8мм', 'Исполнение диода LED': 'трехцветный', 'Цвет диодa LED': 'RGB', 'Длина волны красного цвета " + [char]0x3BB + "d': '620...630нм', 'Длина волны зеленого цвета " + [char]0x3BB + "d': '515...530нм', 'Длина волны голубого цвета': '465...475нм', 'Монтаж': 'THT', 'Характеристики диодов LED': 'программируемый', '#Common #name - #search': 'управляемый', 'Рабочая температура': '-25...80" + [char]0xB0 + "C', 'Кол-во цветов': '16,7M', 'Рабочее напряжение': '1,8...2,2/2,8...3,1/2,9...3,2В'}"
$ws.Range("G15").Value = $g15

# --- G29: add "Напряжение питания" key, drop "Рабочее напряжение" key ---
$g29 = "{'Тип микросхемы': 'микроконтроллер PIC', 'Память программы': '7кБ', 'Объем памяти SRAM': '256Б', 'Объем памяти EEPROM': '256Б', 'Интерфейс': 'MSSP (SPI / I2C)', 'Напряжение питания': '1,8...5,5В DC', 'Монтаж': 'SMD', 'Корпус': 'SO28', 'Кол-во входов/выходов': '25', 'Кол-во таймеров 8бит': '4', 'Кол-во таймеров 16бит': '1', 'Вид архитектуры': 'Harvard  8бит', 'Встроенный генератор': '32МГц', 'Семейство': 'PIC16'}"
$ws.Range("G29").Value = $g29

# --- F20 / F37 / F41: apply hyperlink style + actual hyperlinks ---
$ws.Hyperlinks.Add($ws.Range("F20"), "www.tme.eu/ru/details/m22-fled-rg/panelnye-aksessuary-standartnye/eaton-electric/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F37"), "www.tme.eu/ru/details/ss49e/datchiki-kholla/honeywell/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F41"), "www.tme.eu/ru/details/max660csa+/reguliatory-napriazheniia-mikroskhemy-dc-dc/maxim-integrated/") | Out-Null

# --- Row 49: new product row ---
$ws.Range("A49").Value = "PT10LV10103A2020S"
$ws.Range("B49").Value = "PT10LV-10K"
$ws.Range("C49").Value = "Потенциометр: монтажный; однооборотный,горизонтальный; 10кОм"
$ws.Range("D49").Value = 0.00062
$ws.Range("E49").Value = "static.tme.eu/products_pics/4/8/9/489b0853dca0e1de83ff0b390124fc9c/59967.jpg"
$ws.Range("F49").Value = "www.tme.eu/ru/details/pt10lv-10k/potentsiometry-tht-odnooborotnye/piher/pt10lv10103a2020s/"
$g49 = "{'Тип потенциометра': 'монтажный', 'Вид потенциометра': 'однооборотный', 'Сопротивление': '10кОм', 'Мощность': '150мВт', 'Погрешность': '" + [char]0xB1 + "20%', 'Характеристика': 'линейная', 'Серия потенциометров': 'PT10LV', 'Стандарт потенциометра': '10мм', 'Угол поворота механический': '235 " + [char]0xB1 + "5" + [char]0xB0 + "', 'Крутящий момент': '0,4...2Нсм', 'Рабочее напряжение макс.': '200В', 'Рабочая температура': '-25...70" + [char]0xB0 + "C', 'Монтаж': 'THT', 'Шаг выводов': '5x10мм', 'Материал дорожки': 'углеродистый', 'Угол поворота электрический': '220 " + [char]0xB1 + "20" + [char]0xB0 + "', 'Размеры корпуса': '" + [char]0xD8 + "10,3x4,5мм', 'Характеристики потенциометров': 'шлиц под крестовую отвертку'}"
$ws.Range("G49").Value = $g49
$ws.Range("H49").Value = "www.tme.eu/Document/0c37888a1f172cc56f5b3ed6f5607c24/PIHER_PT10_series.pdf"

# --- sheet view: scroll position & zoom & selection ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 115
$excel.ActiveWindow.ScrollRow = 31
$ws.Range("A49").Select()
